$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 189
$ws.Range("F3").Value = 969
$ws.Range("F4").Value = 1126
$ws.Range("F5").Value = 1567
$ws.Range("F7").Value = 725
$ws.Range("F8").Value = 13050
$ws.Range("F9").Value = 2257
$ws.Range("F11").Value = 295
$ws.Range("F12").Value = 53540
$ws.Range("F13").Value = 1287
$ws.Range("F14").Value = 296
$ws.Range("F15").Value = 298
$ws.Range("F16").Value = 850
$ws.Range("F17").Value = 699
$ws.Range("F18").Value = 351
$ws.Range("F19").Value = 2959
$ws.Range("F20").Value = 840
$ws.Range("F21").Value = 5050
$ws.Range("F22").Value = 1233
$ws.Range("F23").Value = 921
$ws.Range("F25").Value = 32
$ws.Range("F27").Value = 19
$ws.Range("F28").Value = 1174
$ws.Range("F29").Value = 76
$ws.Range("F30").Value = 22
$ws.Range("F31").Value = 139
$ws.Range("F37").Value = 4658
$ws.Range("F39").Value = 4703
$ws.Range("F40").Value = 5664
$ws.Range("F42").Value = 143
$ws.Range("F43").Value = 108
$ws.Range("F44").Value = 196
$ws.Range("F45").Value = 398
$ws.Range("F47").Value = 64
$ws.Range("F48").Value = 4151
$ws.Range("F49").Value = 170

# Sheet 2: 演出
$ws = $wb.Worksheets.Item(2)
$ws.Range("F4").Value = 85
$ws.Range("F12").Value = 1101
$ws.Range("F20").Value = 83

# Sheet 3: 本地生活
$ws = $wb.Worksheets.Item(3)
$ws.Range("F2").Value = 775
$ws.Range("F3").Value = 528
$ws.Range("F4").Value = 137

# Sheet 4: 全部类型
$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 775
$ws.Range("F3").Value = 528
$ws.Range("F4").Value = 189
$ws.Range("F5").Value = 969
$ws.Range("F6").Value = 1126
$ws.Range("F7").Value = 725
$ws.Range("F8").Value = 13050
$ws.Range("F9").Value = 13050
$ws.Range("F10").Value = 2257
$ws.Range("F12").Value = 296
$ws.Range("F13").Value = 850
$ws.Range("F14").Value = 699
$ws.Range("F15").Value = 351
$ws.Range("F16").Value = 2959
$ws.Range("F17").Value = 840
$ws.Range("F18").Value = 85
$ws.Range("F19").Value = 5050
$ws.Range("F20").Value = 1233
$ws.Range("F24").Value = 19
$ws.Range("F26").Value = 1174
$ws.Range("F28").Value = 76
$ws.Range("F29").Value = 22
$ws.Range("F30").Value = 139
$ws.Range("F35").Value = 4658
$ws.Range("F36").Value = 4703
$ws.Range("F37").Value = 5664
$ws.Range("F39").Value = 143
$ws.Range("F40").Value = 196
$ws.Range("F44").Value = 4151
$ws.Range("F46").Value = 83
